# Updates cryptos list price/volume figures (and a few coin-row swaps)
# to match the "Sun Jun 16 18:30:56 UTC 2024" GitHub Actions refresh.
#
# Note: several Price-column values look like plain numbers (e.g. "0.159",
# "6.34"). Excel's Range.Value setter auto-converts such text to a numeric
# type, which would change the cell's stored kind from text to number.
# To keep them as text (matching the source data, which uses '.' as a
# thousands separator in some rows, e.g. "66.623.80"), those values are
# written with a leading apostrophe (forces text) and the cell's original
# Style is restored immediately after so no stray "quote prefix" / text
# number-format style gets left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.623.80'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '3.593.33'
$ws.Range('E3').Value = '  +0.82%  '
$s = $ws.Range('D5').Style
$ws.Range('D5').Value = '''609.53'
$ws.Range('D5').Style = $s
$ws.Range('E5').Value = '  +0.39%  '
$s = $ws.Range('D6').Style
$ws.Range('D6').Value = '''148.25'
$ws.Range('D6').Style = $s
$ws.Range('E6').Value = '  +2.32%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.47%  '
$s = $ws.Range('D9').Style
$ws.Range('D9').Value = '''8.06'
$ws.Range('D9').Style = $s
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = '4.200.79'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('E13').Value = '  +0.68%  '
$s = $ws.Range('D14').Style
$ws.Range('D14').Value = '''29.90'
$ws.Range('D14').Style = $s
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').Value = '3.620.21'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').Value = '66.691.55'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('E17').Value = '  +0.77%  '
$s = $ws.Range('D18').Style
$ws.Range('D18').Value = '''11.45'
$ws.Range('D18').Style = $s
$ws.Range('E18').Value = '  -0.10%  '
$s = $ws.Range('D19').Style
$ws.Range('D19').Value = '''6.34'
$ws.Range('D19').Style = $s
$ws.Range('E19').Value = '  +2.05%  '
$s = $ws.Range('D20').Style
$ws.Range('D20').Value = '''15.11'
$ws.Range('D20').Style = $s
$ws.Range('E20').Value = '  +1.36%  '
$s = $ws.Range('D21').Style
$ws.Range('D21').Value = '''427.77'
$ws.Range('D21').Style = $s
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('E22').Value = '  +1.51%  '
$s = $ws.Range('D23').Style
$ws.Range('D23').Value = '''78.86'
$ws.Range('D23').Style = $s
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '3.732.34'
$ws.Range('E25').Value = '  +0.04%  '
$s = $ws.Range('D26').Style
$ws.Range('D26').Value = '''0.0000122'
$ws.Range('D26').Style = $s
$ws.Range('E26').Value = '  +2.69%  '
$s = $ws.Range('D27').Style
$ws.Range('D27').Value = '''8.30'
$ws.Range('D27').Style = $s
$ws.Range('E27').Value = '  +3.32%  '
$s = $ws.Range('D28').Style
$ws.Range('D28').Value = '''9.33'
$ws.Range('D28').Style = $s
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$s = $ws.Range('D31').Style
$ws.Range('D31').Value = '''0.159'
$ws.Range('D31').Style = $s
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('B32').Value = 'RenzoRestakedETH'
$ws.Range('C32').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D32').Value = '3.589.61'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$s = $ws.Range('D33').Style
$ws.Range('D33').Value = '''1.46'
$ws.Range('D33').Style = $s
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$s = $ws.Range('D34').Style
$ws.Range('D34').Value = '''25.49'
$ws.Range('D34').Style = $s
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$s = $ws.Range('D37').Style
$ws.Range('D37').Value = '''5.64'
$ws.Range('D37').Style = $s
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$s = $ws.Range('D38').Style
$ws.Range('D38').Value = '''1.71'
$ws.Range('D38').Style = $s
$ws.Range('E38').Value = '  -2.85%  '
$s = $ws.Range('D39').Style
$ws.Range('D39').Value = '''177.97'
$ws.Range('D39').Style = $s
$ws.Range('E39').Value = '  +4.50%  '
$ws.Range('E40').Value = '  +0.28%  '
$s = $ws.Range('D41').Style
$ws.Range('D41').Value = '''5.23'
$ws.Range('D41').Style = $s
$ws.Range('E41').Value = '  +0.41%  '
$s = $ws.Range('D42').Style
$ws.Range('D42').Value = '''0.899'
$ws.Range('D42').Style = $s
$ws.Range('E42').Value = '  +0.16%  '
$s = $ws.Range('D43').Style
$ws.Range('D43').Value = '''1.92'
$ws.Range('D43').Style = $s
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('E44').Value = '  +8.22%  '
$s = $ws.Range('D45').Style
$ws.Range('D45').Value = '''1.00'
$ws.Range('D45').Style = $s
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -1.80%  '
$s = $ws.Range('D47').Style
$ws.Range('D47').Value = '''24.34'
$ws.Range('D47').Style = $s
$s = $ws.Range('D48').Style
$ws.Range('D48').Value = '''25.03'
$ws.Range('D48').Style = $s
$ws.Range('E48').Value = '  -3.24%  '
$s = $ws.Range('D49').Style
$ws.Range('D49').Value = '''7.18'
$ws.Range('D49').Style = $s
$ws.Range('E49').Value = '  +0.56%  '
$s = $ws.Range('D50').Style
$ws.Range('D50').Value = '''0.951'
$ws.Range('D50').Style = $s
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('E51').Value = '  -1.54%  '
